# CompStat weekly refresh: new report week + updated crime-complaint figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: some cells in the "% Chg" columns display a literal dash/placeholder
# text ("0") instead of a real number when there's no meaningful comparison.
# Typing a numeric-looking string via .Value normally gets auto-coerced back
# into a number (like real Excel), so first flip the cell to Text format,
# write the literal digit, then restore the original General-formatted
# "label" style (copied from a neighboring cell that already holds that same
# style) without disturbing the freshly-written text value.
# ---------------------------------------------------------------------------
function Set-DashText($cellRef) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = "0"
    $ws.Range("C15").Copy()
    $ws.Range($cellRef).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $excel.CutCopyMode = $false
}

# --- Masthead: volume/issue number and reporting week -----------------------
$ws.Range("A8").Value = "Volume 30   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/11/2023  Through  12/17/2023"

# --- Murder ------------------------------------------------------------------
Set-DashText "C14"

# --- Rape ---------------------------------------------------------------------
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0

# --- Robbery -------------------------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -83.333333333333
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -57.142857142857
$ws.Range("I16").Value = 137
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 14.166666666666
$ws.Range("L16").Value = 29.245283018867
$ws.Range("M16").Value = -18.452380952381
$ws.Range("N16").Value = -81.284153005464

# --- Fel. Assault ----------------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 40
$ws.Range("I17").Value = 218
$ws.Range("J17").Value = 176
$ws.Range("K17").Value = 23.863636363636
$ws.Range("L17").Value = 27.485380116959
$ws.Range("M17").Value = 66.412213740458
$ws.Range("N17").Value = -26.599326599326

# --- Burglary ----------------------------------------------------------------
Set-DashText "C18"
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -50
$ws.Range("J18").Value = 182
$ws.Range("K18").Value = 11.538461538461
$ws.Range("L18").Value = 41.958041958042
$ws.Range("M18").Value = -29.757785467128
$ws.Range("N18").Value = -87.359900373599

# --- Gr. Larceny ---------------------------------------------------------------
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 75
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 2.380952380952
$ws.Range("I19").Value = 621
$ws.Range("J19").Value = 702
$ws.Range("K19").Value = -11.538461538461
$ws.Range("L19").Value = 7.8125
$ws.Range("M19").Value = 42.105263157894
$ws.Range("N19").Value = -17.857142857142

# --- G.L.A. (note: D20/E20 flip from placeholder text to real numbers) -------
$ws.Range("C20").Value = 3
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 3
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 75
$ws.Range("I20").Value = 159
$ws.Range("J20").Value = 143
$ws.Range("K20").Value = 11.188811188811
$ws.Range("L20").Value = 41.964285714285
$ws.Range("M20").Value = 4.605263157894
$ws.Range("N20").Value = -90.991501416430

# --- TOTAL ---------------------------------------------------------------------
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 84
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = -1.176470588235
$ws.Range("I21").Value = 1360
$ws.Range("J21").Value = 1345
$ws.Range("K21").Value = 1.115241635687
$ws.Range("L21").Value = 20.353982300885
$ws.Range("M21").Value = 13.807531380753
$ws.Range("N21").Value = -73.790711119676

# --- Transit (C22 flips from placeholder text to a real number) --------------
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 19
$ws.Range("K22").Value = 58.333333333333
$ws.Range("L22").Value = 11.764705882352
$ws.Range("M22").Value = -29.629629629629

# --- Petit Larceny -------------------------------------------------------------
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 118
$ws.Range("G24").Value = 142
$ws.Range("H24").Value = -16.901408450704
$ws.Range("I24").Value = 1635
$ws.Range("J24").Value = 1824
$ws.Range("K24").Value = -10.361842105263
$ws.Range("L24").Value = 32.281553398058
$ws.Range("M24").Value = 73.015873015873

# --- Misd. Assault ----------------------------------------------------------
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 66.666666666666
$ws.Range("F25").Value = 53
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 55.882352941176
$ws.Range("I25").Value = 610
$ws.Range("J25").Value = 449
$ws.Range("K25").Value = 35.857461024498
$ws.Range("L25").Value = 40.877598152424
$ws.Range("M25").Value = 44.208037825059

# --- UCR Rape* -----------------------------------------------------------------
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = -50

# --- Other Sex Crimes --------------------------------------------------------
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 49
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = -2
$ws.Range("L27").Value = -19.672131147541

# --- Shooting Vic. / Shooting Inc. --------------------------------------------
$ws.Range("M28").Value = 0
$ws.Range("M29").Value = 0
